$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so number-like price strings
# (e.g. "0.9995", "238.27") are not auto-converted to numeric values,
# matching the source data which stores prices as inline strings.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '25.877.29'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '1.740.45'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '238.27'
$ws.Range("E5").Value = '  +3.49%  '
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").Value = '0.5159'
$ws.Range("E7").Value = '  -1.23%  '
$ws.Range("D8").Value = '0.2735'
$ws.Range("E8").Value = '  -0.49%  '
$ws.Range("D9").Value = '0.06140'
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("D10").Value = '1.737.98'
$ws.Range("E10").Value = '  +0.12%  '
$ws.Range("D11").Value = '0.07167'
$ws.Range("E11").Value = '  +1.89%  '
$ws.Range("D12").Value = '0.6447'
$ws.Range("E12").Value = '  +1.48%  '
$ws.Range("D13").Value = '14.93'
$ws.Range("E13").Value = '  -0.48%  '
$ws.Range("D14").Value = '4.591'
$ws.Range("E14").Value = '  +1.43%  '
$ws.Range("D15").Value = '77.28'
$ws.Range("E15").Value = '  +0.81%  '
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").Value = '0.9994'
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("D18").Value = '25.893.11'
$ws.Range("E18").Value = '  +0.29%  '
$ws.Range("E19").Value = '  +2.13%  '
$ws.Range("D20").Value = '0.000006768'
$ws.Range("E20").Value = '  +1.79%  '
$ws.Range("D21").Value = '1.962.50'
$ws.Range("E21").Value = '  +0.34%  '
$ws.Range("D22").Value = '4.264'
$ws.Range("E22").Value = '  +1.86%  '
$ws.Range("D23").Value = '8.666'
$ws.Range("E23").Value = '  -1.00%  '
$ws.Range("D24").Value = '5.233'
$ws.Range("E24").Value = '  +1.45%  '
$ws.Range("D25").Value = '138.35'
$ws.Range("E25").Value = '  -0.77%  '
$ws.Range("D26").Value = '1.508'
$ws.Range("E26").Value = '  +0.36%  '
$ws.Range("D27").Value = '15.12'
$ws.Range("E27").Value = '  +0.78%  '
$ws.Range("D28").Value = '1.762'
$ws.Range("E28").Value = '  -1.01%  '
$ws.Range("D29").Value = '105.81'
$ws.Range("E29").Value = '  +3.76%  '
$ws.Range("D30").Value = '3.954'
$ws.Range("E30").Value = '  +6.67%  '
$ws.Range("D31").Value = '0.08292'
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("D32").Value = '3.638'
$ws.Range("E32").Value = '  +3.96%  '
$ws.Range("D33").Value = '0.04583'
$ws.Range("E33").Value = '  +2.67%  '
$ws.Range("D34").Value = '2.663'
$ws.Range("E34").Value = '  +2.24%  '
$ws.Range("D35").Value = '0.9882'
$ws.Range("E35").Value = '  +1.67%  '
$ws.Range("D36").Value = '0.6178'
$ws.Range("E36").Value = '  +0.28%  '
$ws.Range("D37").Value = '2.686'
$ws.Range("D38").Value = '0.01612'
$ws.Range("E38").Value = '  +2.76%  '
$ws.Range("D39").Value = '1.929'
$ws.Range("E39").Value = '  +1.51%  '
$ws.Range("D40").Value = '0.9992'
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("D41").Value = '97.74'
$ws.Range("E41").Value = '  -2.31%  '
$ws.Range("D42").Value = '0.3829'
$ws.Range("E42").Value = '  +0.12%  '
$ws.Range("D43").Value = '0.7385'
$ws.Range("E43").Value = '  +2.42%  '
$ws.Range("D44").Value = '4.961'
$ws.Range("E44").Value = '  -0.94%  '
$ws.Range("E45").Value = '  -0.15%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '6.196'
$ws.Range("E46").Value = '  +0.44%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '0.05259'
$ws.Range("E47").Value = '  -1.61%  '
$ws.Range("D48").Value = '54.84'
$ws.Range("E48").Value = '  +3.27%  '
$ws.Range("D49").Value = '30.48'
$ws.Range("E49").Value = '  +1.81%  '
$ws.Range("D50").Value = '7.572'
$ws.Range("E50").Value = '  -0.24%  '
$ws.Range("D51").Value = '0.3405'
$ws.Range("E51").Value = '  +0.92%  '

# Restore the default cell style so no stray formatting is left behind
# on the price column (keeps the cell style index the same as before).
$priceRange.Style = "Normal"
